$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fix capitalization of names (column B, rows 2-5) and fix "Mumabai" typo in row 5's address
$ws.Range("B2").Value = "Sachin"
$ws.Range("B3").Value = "Satvik"
$ws.Range("B4").Value = "Sania"
$ws.Range("B5").Value = "Dishant"

$ws.Range("C5").Value = "802/ Gunjan nagar/ Andheri , Mumbai"

# Auto-fit all used columns so widths match "bestFit" sizing
$ws.Range("A1:D5").EntireColumn.AutoFit() | Out-Null

# Update the active selection to D11
$ws.Range("D11").Select()
